# Script: reorders three existing match rows (3,4,5) and appends a new
# match row (80) to the Croatia Prva NL 2023-2024 results sheet, matching
# the "Atualizado por script em 10-11-2023 14:45" refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: cyclic rotation of the match data (columns F:V) held in
# rows 3, 4 and 5. The row identifiers/date in A:E stay untouched; only
# the match info moves:
#   row3 <- old row5, row4 <- old row3, row5 <- old row4
$old3 = $ws.Range("F3:V3").Value()
$old4 = $ws.Range("F4:V4").Value()
$old5 = $ws.Range("F5:V5").Value()

$ws.Range("F3:V3").Value = $old5
$ws.Range("F4:V4").Value = $old3
$ws.Range("F5:V5").Value = $old4

# --- Step 2: append the new match as row 80.
$ws.Range("A80").Value = 79
$ws.Range("B80").Value = "croatia"
$ws.Range("C80").Value = "prva-nl"
$ws.Range("D80").Value = "2023-2024"
$ws.Range("E80").Value = 45240.58333333334
$ws.Range("F80").Value = "Solin"
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = "Zrinski Jurjevac"
$ws.Range("I80").Value = 2
$ws.Range("J80").Value = 2.47
$ws.Range("K80").Value = "09/11/2023 02:12"
$ws.Range("L80").Value = 3.18
$ws.Range("M80").Value = "10/11/2023 13:20"
$ws.Range("N80").Value = 3.08
$ws.Range("O80").Value = "09/11/2023 02:12"
$ws.Range("P80").Value = 3
$ws.Range("Q80").Value = "10/11/2023 13:20"
$ws.Range("R80").Value = 2.6
$ws.Range("S80").Value = "09/11/2023 02:12"
$ws.Range("T80").Value = 2.37
$ws.Range("U80").Value = "10/11/2023 13:20"
$ws.Range("V80").Value = "https://www.betexplorer.com/football/croatia/prva-nl/solin-zrinski-jurjevac/EeKzxxfL/"

# Match the existing formatting used by the other data rows: column A is
# bold/bordered/centered (same style as the other "Indice" cells) and
# column E keeps the custom date-time number format.
$ws.Range("A2").Copy()
$ws.Range("A80").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E80").PasteSpecial(-4122)

$excel.CutCopyMode = 0
